# Adding parameters from excel file part2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "Result" values for the existing rows (2-6) + two inserted G cells.
# The e-mail text in column D is updated too; the pre-existing mailto
# hyperlinks (and their relationship ids/targets) are left untouched.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "Pass - That user exists"

$ws.Range("D3").Value = "mando200720252@gmail.com"
$ws.Range("G3").Value = "Pass - User Created"

$ws.Range("D4").Value = "dldg20072025@gmail.com"
$ws.Range("G4").Value = "Pass - First  and last Name are not valid"

$ws.Range("D5").Value = "mando270420252@gmail.com"
$ws.Range("G5").Value = "Pass - First Name is not valid!"

$ws.Range("D6").Value = "mando270420252@gmail.com"
$ws.Range("G6").Value = "Pass - Last Name is not valid!!"

# ---------------------------------------------------------------------
# Row 7 - test_negative_incorrect_format_email
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "test_negative_incorrect_format_email"
$ws.Range("B7").Value = "Pedro"
$ws.Range("C7").Value = "Pascal"
$ws.Range("D7").Value = "dldg210461gmail.com"
$ws.Range("E7").Value = "Groguforever123"
$ws.Range("F7").Value = "Groguforever123"
$ws.Range("G7").Value = "Pass - email is not valid!!"

$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:dldg210461@gmail.com", "", "", "dldg210461@gmail.com") | Out-Null
$ws.Range("D7").Style = $ws.Range("D2").Style

# ---------------------------------------------------------------------
# Row 8 - test_negative_password_Strength
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "test_negative_password_Strength"
$ws.Range("D8").Style = $ws.Range("D2").Style
$ws.Range("E8").Value = "Grogu"
$ws.Range("G8").Value = "Pass - Not accurate Password!"

# ---------------------------------------------------------------------
# Row 9 - test_positive_password_Strength (new monospaced font style)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "test_positive_password_Strength"
$f9 = $ws.Range("A9").Font
$f9.Name = "JetBrains Mono"
$f9.Family = 3
$f9.Size = 10
$f9.Color = 16099414
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("E9").Value = "Groguforever123"
$ws.Range("G9").Value = "Pass - Accurate Password!"

# ---------------------------------------------------------------------
# Row 10 - test_negative_password_confirm_diff
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "test_negative_password_confirm_diff"
$ws.Range("B10").Value = "Pedro"
$ws.Range("C10").Value = "Pascal"
$ws.Range("D10").Value = "mando20042028@gmail.com"
$ws.Range("E10").Value = "Groguforever123"
$ws.Range("F10").Value = "Groguforever123*"
$ws.Range("G10").Value = "Pass - Password and confirm are different"

$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:mando20042028@gmail.com") | Out-Null
$ws.Range("D10").Style = $ws.Range("D2").Style

# ---------------------------------------------------------------------
# Move the active selection to D12, matching the saved view state
# ---------------------------------------------------------------------
$ws.Range("D12").Select() | Out-Null
